# Deep Audit: Remove all hardcoded defaults across the workbook
# Sets the previously hardcoded default numeric inputs to 0 on each sheet.

$wb = $excel.ActiveWorkbook

# --- WORKFORCE_PLANNING sheet ---
$wsWorkforce = $wb.Worksheets.Item("WORKFORCE_PLANNING")
$wsWorkforce.Range("B5").Value = 0   # Est. Hiring Fee (per worker)
$wsWorkforce.Range("B6").Value = 0   # Est. Severance (per worker)
$wsWorkforce.Range("D10").Value = 0  # Est. Turnover % - Center
$wsWorkforce.Range("D11").Value = 0  # Est. Turnover % - West
$wsWorkforce.Range("D12").Value = 0  # Est. Turnover % - North
$wsWorkforce.Range("D13").Value = 0  # Est. Turnover % - East
$wsWorkforce.Range("D14").Value = 0  # Est. Turnover % - South

# --- COMPENSATION_STRATEGY sheet ---
$wsComp = $wb.Worksheets.Item("COMPENSATION_STRATEGY")
$wsComp.Range("B6").Value = 0    # Inflation Rate %
$wsComp.Range("B7").Value = 0    # Target Purchasing Power Increase %

$wsComp.Range("B11").Value = 0   # Previous Salary - Center
$wsComp.Range("D11").Value = 0   # Proposed New Salary - Center
$wsComp.Range("B12").Value = 0   # Previous Salary - West
$wsComp.Range("D12").Value = 0   # Proposed New Salary - West
$wsComp.Range("B13").Value = 0   # Previous Salary - North
$wsComp.Range("D13").Value = 0   # Proposed New Salary - North
$wsComp.Range("B14").Value = 0   # Previous Salary - East
$wsComp.Range("D14").Value = 0   # Proposed New Salary - East
$wsComp.Range("B15").Value = 0   # Previous Salary - South
$wsComp.Range("D15").Value = 0   # Proposed New Salary - South

$wsComp.Range("B23").Value = 0   # Training Budget (% of Payroll)
$wsComp.Range("B24").Value = 0   # Health Insurance (% of Payroll)
$wsComp.Range("B25").Value = 0   # Profit Sharing (% of Net Profit)
$wsComp.Range("B26").Value = 0   # Personal Days (per Worker)
$wsComp.Range("B27").Value = 0   # Union Representatives

# --- LABOR_COST_ANALYSIS sheet ---
$wsLabor = $wb.Worksheets.Item("LABOR_COST_ANALYSIS")
$wsLabor.Range("B4").Value = 0   # INPUT: Estimated Net Profit (for Profit Sharing)
